$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Leonardo Parisi "
$ws.Range("B12").Value = "Stefano Tita | Clitoriders"
$ws.Range("C12").Value = "Federico  Mortillaro | Clitoriders"
$ws.Range("D12").Value = "Federico  Manica | iMontagna"
$ws.Range("E12").Value = "Nicholas Marzadro | SBARX"
$ws.Range("F12").Value = "Maverick  Bertolini | A.C. Denti"
